$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.535700000000005
$ws.Range("B21").Value = 9.573200000000003
$ws.Range("B23").Value = 9.035400000000006
$ws.Range("B25").Value = 5.926799999999999
$ws.Range("E27").Value = 16.73229999999998
$ws.Range("E31").Value = 16.5218
$ws.Range("E39").Value = 16.06409999999999
$ws.Range("E48").Value = 17.50830000000001
$ws.Range("E51").Value = 17.2364
$ws.Range("E52").Value = 17.0353
$ws.Range("B53").Value = 5.513599999999998
$ws.Range("E55").Value = 16.61960000000001
$ws.Range("E56").Value = 16.2625
$ws.Range("B57").Value = 4.888399999999997
$ws.Range("E57").Value = 16.67090000000001
$ws.Range("B59").Value = 4.777099999999998
$ws.Range("B69").Value = 5.282799999999997
$ws.Range("E73").Value = 17.15530000000001
$ws.Range("B79").Value = 8.936000000000002
$ws.Range("B83").Value = 5.226799999999998
$ws.Range("E89").Value = 17.29190000000002
$ws.Range("E90").Value = 16.65609999999998
$ws.Range("B93").Value = 5.487699999999998
